$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the 100 new videoID rows (rows 252-351) that were scraped into the sheet.
$newIds = @(
    "JAjxq93njvs",
    "Bi0taCrL22Q",
    "Vkack7OKYjQ",
    "o7FmGwoDJvc",
    "jJUvn-i73zI",
    "1AMFwO7BGW8",
    "8DkcTriimSM",
    "G5iGJHVtFqs",
    "Mh5dKiiZRlg",
    "3he4lvoghMg",
    "xTsREJ6KEfg",
    "kaqxOKmq0ZI",
    "muT3_bTXR_w",
    "5tQPftbhCVo",
    "tivtBueDg_c",
    "k7PmttZtFU0",
    "m5zhPegPUc8",
    "Su_SV29obhE",
    "UmTBOw5ITeU",
    "Er2IIYpGgS0",
    "ouPsgmEKyww",
    "19IdcJU8524",
    "G5_YyZUOQL8",
    "jCIf07qU_mQ",
    "xhIFmxslJCQ",
    "L3ibkV6t0NM",
    "4vfIC0dZ8jM",
    "MrNcHCki8mE",
    "YDdN4vm_59c",
    "u6Q8ywGzkAI",
    "DdPWX4RqDOY",
    "YSrjyXqA5cM",
    "jn73kYaLmKM",
    "2HS3sFyywnM",
    "O6mZnetxUX8",
    "vyB_2BnouRs",
    "x-axYKSh-jc",
    "2nUGLtfiTbk",
    "qb6UuECcPV0",
    "YpUgnlZkPEI",
    "WMSpjDc8zH8",
    "sxn5n89S7Xw",
    "VKJZ4PCWgY8",
    "Bd2fE3YO7Aw",
    "4SL11QJZzCI",
    "0ONPlWHl9Gk",
    "L-PkACES1qc",
    "o8hO9R_Xs7Q",
    "Qs0KZwhKdwo",
    "moSE2_ur994",
    "BlO34GiRoJ0",
    "JM1S2ZkNVrQ",
    "L4B7z8dHDRc",
    "cKW-pp1pprg",
    "sWVPnyyLnFw",
    "DPDUDPCttfc",
    "OGiBcHF-wOU",
    "VkqiSzoVtUA",
    "Id1rwD0RU5k",
    "AHIlxaPZd5w",
    "mxWIy8q_vQw",
    "aslFiO7vThQ",
    "ro5arbsZfMs",
    "goCGfcCjtSo",
    "MAm0RLQpYas",
    "qKKNbYv3LrM",
    "vO4QeHR9SXM",
    "w-253j0BjLk",
    "hddRxBhGiMI",
    "YqKYpgZ9FWU",
    "gTo-lPOGPdg",
    "PzHZGefMx9Q",
    "btFETPVpDcw",
    "nSNkwIvlWJU",
    "LWHMtkbCxd0",
    "YJOmxjtG9IM",
    "jEWaLRi_G60",
    "GlUp3MCIQIk",
    "OqGH_6mLQjI",
    "S1ABZywPONY",
    "eO6e0b4i93U",
    "1CLhjDLUl4c",
    "xDi_TTbk3-c",
    "Z0L0xWZhYTE",
    "t8rMgGSgHv0",
    "sJj3q6njVbM",
    "ImO6sDdYB8U",
    "3Bi-0s8djNE",
    "04kGIbHUIho",
    "Y1E6R13YD1Y",
    "7hTgjEzl4Bk",
    "BW_MAa5L9lg",
    "16RX6rKJle8",
    "hN0xpuih8zo",
    "mcTpep5uoZ8",
    "uqv-oC-Izi8",
    "_bhXh1KxVHc",
    "nSI_2yvBAIU",
    "8a5Z1Idf-Jo",
    "dkf9LNwo2kM"
)

$startRow = 252
for ($i = 0; $i -lt $newIds.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $newIds[$i]
}

# Update the view so it reflects where the user ended up after pasting the new data.
[void]$ws.Range("E311").Select()
$excel.ActiveWindow.ScrollRow = 306
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 115
